$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new registry row (row 51) for GTSm / ISO 21812
$ws.Range("A51").Value = "GTSm"
$ws.Range("B51").Value = "Debbie"
$ws.Range("C51").Value = "Orf"

$ws.Range("D51").Value = "dorf@aptech.com"
$ws.Hyperlinks.Add($ws.Range("D51"), "mailto:dorf@aptech.com")
$ws.Range("D51").Style = "Hyperlink"

$ws.Range("E51").Value = "ISO TC130/WG2 as described in ISO 21812"
$ws.Range("F51").Value = 43616

# Update the view's active selection to the next empty row, as in the source workbook
$null = $ws.Range("A52").Select()
